$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: build a pkg:package wrapper around a body fragment and InsertXML it
# into the supplied Range. Word's InsertXML call only behaves correctly (i.e.
# replaces exactly the addressed Range) when that Range's end coincides with
# the end of its containing paragraph's text (just before the paragraph
# mark); otherwise the new content is appended at the paragraph's end while
# the addressed text is simply deleted. To stay safe we therefore always
# extend the target Range through the end of the paragraph and re-emit any
# untouched trailing run content verbatim as part of the replacement XML.
# ---------------------------------------------------------------------------
function New-PkgXml([string]$bodyFragment) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$bodyFragment
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ===========================================================================
# Change 1: Blood Test Manager / INR app paragraph - rewrite the description
# of the tech stack, ending on "MongoDB Atlas." which is also the very last
# text of the paragraph, so a straightforward Range replace is safe here.
# ===========================================================================

$afterFullStack = $d.Content
$afterFullStack.Find.Execute("full-stack", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFullStack.Collapse(0)
$startPos = $afterFullStack.Start

$endRange = $d.Content
$endRange.Find.Execute("MongoDB Atlas.", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endRange.End

$target1 = $d.Range($startPos, $endPos)

$body1 = @"
<w:p>
<w:r><w:t xml:space="preserve"> iOS application to manage and visualize my INR (International Normalized Ratio) blood test results. </w:t></w:r>
<w:r><w:t xml:space="preserve">The first iteration of this project utilized Storyboard UI, Node.js server, and MongoDB database. The second iteration utilized </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>SwiftUI</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> and Firebase to implement multiple user accounts and authentication logic.</w:t></w:r>
</w:p>
"@

$target1.InsertXML((New-PkgXml $body1))

# ===========================================================================
# Change 2: Service Ticket Database paragraph - split the long run that
# contains "... Mean Time " into two runs, inserting a lastRenderedPageBreak
# marker at the start of the second one. The edit point is NOT at the end of
# the paragraph, so we must extend our replacement Range through the rest of
# the paragraph (the untouched "to" / " Repair (MTTR)." runs) and re-emit
# that trailing content verbatim (including the original w:rsidR attribute)
# so nothing is lost.
# ===========================================================================

$afterLocal = $d.Content
$afterLocal.Find.Execute("local", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterLocal.Collapse(0)
$startPos2 = $afterLocal.Start

$paraRange = $d.Content
$paraRange.Find.Execute("Created a Java program", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraRange.Expand(4) | Out-Null
$paraEnd2 = $paraRange.End - 1

$target2 = $d.Range($startPos2, $paraEnd2)

$body2 = @"
<w:p>
<w:r><w:t xml:space="preserve"> MySQL database with over 15,000 randomly generated service tickets. This database was then </w:t></w:r>
<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">connected to a business intelligence software, Looker Studio, to visualize and derive various metrics from the database including ticket Priority and Mean Time </w:t></w:r>
<w:r w:rsidR="00001027"><w:t>to</w:t></w:r>
<w:r><w:t xml:space="preserve"> Repair (MTTR).</w:t></w:r>
</w:p>
"@

$target2.InsertXML((New-PkgXml $body2))

Write-Output "Done."
